# Integre el registro de proyectos tanto individual como masivos
# Insert two new header columns ("Correo del Coordinador" and
# "Teléfono del Coordinador") before the existing "Carrera del Coordinador"
# column, and adjust the view state to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former H1 header ("Carrera del Coordinador") moves to J1, and the two
# new headers take its old spot (H1) plus the new I1 slot.
$ws.Range("J1").Value = "Carrera del Coordinador"
$ws.Range("H1").Value = "Correo del Coordinador"
$ws.Range("I1").Value = "Teléfono del Coordinador"

# New columns need the same kind of explicit width the other header columns
# already carry.
$ws.Columns.Item(9).ColumnWidth = 25.666666666666668
$ws.Columns.Item(10).ColumnWidth = 22.166666666666668

# Update the view: zoomed in, selection resting on D7.
$excel.ActiveWindow.Zoom = 145
$ws.Range("D7").Select() | Out-Null
